# Add the list of student IDs (MSSV) below the existing header/sample rows,
# then format them: Times New Roman font, centered + wrapped text, thin
# black border all around, and auto-fit the column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids = @(
    23030112,23030159,23030101,23030154,23030104,23030113,23030072,23030130,
    23030247,23030100,23030122,23030244,23030220,23030230,23030073,23030034,
    23030002,23030025,23030084,23030118,23030083,23030075,23030026,23030148,
    23030139,23030087,23030152,23030114,23030116,23030115,23030169,23030176,
    23030180,23030125,23030186,23030153
)

$startRow = 5
$row = $startRow
foreach ($id in $ids) {
    $ws.Cells.Item($row, 1).Value = $id
    $row = $row + 1
}
$endRow = $row - 1

$dataRange = $ws.Range("A$startRow`:A$endRow")

# Font: Times New Roman, size 11
$dataRange.Font.Name = "Times New Roman"
$dataRange.Font.Size = 11

# Alignment: centered horizontally and vertically, wrap text
$dataRange.HorizontalAlignment = -4108  # xlCenter
$dataRange.VerticalAlignment = -4108    # xlCenter
$dataRange.WrapText = $true

# Thin black border around every cell
$dataRange.Borders.Color = 0

# Auto-fit column A to the new content
$ws.Columns.Item(1).AutoFit()

Write-Host "Added $($ids.Count) student IDs to A$startRow`:A$endRow"
